# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.767.50'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '3.830.36'
$ws.Range("E3").Value = '  +2.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("D7").Value = '3.832.24'
$ws.Range("E7").Value = '  +2.24%  '
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.47'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.41%  '
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.83'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.13%  '
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("D15").Value = '4.463.86'
$ws.Range("E15").Value = '  +2.35%  '
$ws.Range("D16").Value = '3.825.88'
$ws.Range("E16").Value = '  +2.11%  '
$ws.Range("D17").Value = '69.848.78'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("E19").Value = '  -3.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '507.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.59%  '
$ws.Range("E24").Value = '  -2.79%  '
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("E26").Value = '  +4.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.51'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.74'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("E34").Value = '  -1.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  -1.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.12'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.46%  '
$ws.Range("E38").Value = '  +5.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '485.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +14.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.336'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.05'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.49%  '
$ws.Range("E42").Value = '  -2.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = '2.934.66'
$ws.Range("E46").Value = '  -2.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0361'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '139.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '27.05'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.44'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.52%  '
